$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 13 (shifts existing rows 13-46 down to 14-47,
# matching the diff's row-shift pattern and growing the used range to A1:T47).
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly price record
# (Chirimoya, Primera, Macroferia Regional de Talca - 2021-10-26).
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Macroferia Regional de Talca"
$ws.Range("C13").Value = "Maule"
$ws.Range("D13").Value = 44495
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = "Fruta"
$ws.Range("G13").Value = 100107
$ws.Range("H13").Value = "Otros"
$ws.Range("I13").Value = 100107002
$ws.Range("J13").Value = "Chirimoya"
$ws.Range("K13").Value = "Cultivar IV Región"
$ws.Range("L13").Value = "Primera"
$ws.Range("M13").Value = 150
$ws.Range("N13").Value = 25000
$ws.Range("O13").Value = 25000
$ws.Range("P13").Value = 25000
$ws.Range("Q13").Value = "`$/bandeja 10 kilos"
$ws.Range("R13").Value = "Provincia de Limarí"
$ws.Range("S13").Value = 2500
$ws.Range("T13").Value = 10
